$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 14 (currently the single 2017-2022 Northwestern entry) ---
# Old: 2017 | 2022 | Eratus Otis Haven Professor | - | - | Northwestern University
# New: 2017 | present | Professor | by courtesy | Dept. of Molecular Biosciences | Northwestern University
$ws.Range("C14").Value = "present"
$ws.Range("D14").Value = "Professor"
$ws.Range("E14").Value = "by courtesy"
$ws.Range("F14").Value = "Dept. of Molecular Biosciences"
$ws.Range("G14").Value = "Northwestern University"

# --- Insert a new row 15 for the (corrected) named professorship, split out of old row 14 ---
# New: 2017 | 2022 | Erastus Otis Haven Professor | - | - | Northwestern University
$ws.Range("A15").Formula = "=A14+1"
$ws.Range("B15").Value = 2017
$ws.Range("C15").Value = 2022
$ws.Range("D15").Value = "Erastus Otis Haven Professor"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "Northwestern University"

# Match the numeric/index style used by the rest of column A
$ws.Range("A14").Copy($ws.Range("A15"))
$ws.Range("A15").Formula = "=A14+1"

# Update selection to mirror the new active cell after editing
[void]$ws.Range("D16").Select()
